$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two middle "Periodo Mora" rows (periods 2506 and 2505).
# This shifts the old row 19 (period 2504, which carries the table's
# bottom-border styling) up to become the new row 17.
$ws.Range("17:18").Delete()

# Row 16 now shows the first remaining period: 2506
$ws.Range("E16").Value = "2506"

# Row 17 (previously row 19 / period 2504) becomes the second, last
# period row: 2507, with the regular (non-reduced) amount.
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 160000

# Update the account summary figures for the new, smaller period set.
$ws.Range("E11").Value = 320000
$ws.Range("F13").Value = 2
